$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1575
$ws1.Range("F6").Value = 762
$ws1.Range("F7").Value = 678
$ws1.Range("F8").Value = 1278
$ws1.Range("F9").Value = 2520
$ws1.Range("F10").Value = 1327
$ws1.Range("F12").Value = 2293
$ws1.Range("F14").Value = 704
$ws1.Range("F15").Value = 6162
$ws1.Range("F23").Value = 2036
$ws1.Range("F25").Value = 667
$ws1.Range("F28").Value = 5197
$ws1.Range("F30").Value = 1232
$ws1.Range("F32").Value = 3648
$ws1.Range("F34").Value = 1155
$ws1.Range("F38").Value = 953
$ws1.Range("F39").Value = 367
$ws1.Range("F44").Value = 879
$ws1.Range("F45").Value = 1039

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 435
$ws2.Range("F11").Value = 374
$ws2.Range("F23").Value = 349
$ws2.Range("F25").Value = 167
$ws2.Range("F35").Value = 37

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F13").Value = 229
$ws3.Range("F14").Value = 1106

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 1575
$ws4.Range("F9").Value = 762
$ws4.Range("F11").Value = 678
$ws4.Range("F12").Value = 2520
$ws4.Range("F14").Value = 1327
$ws4.Range("F16").Value = 2293
$ws4.Range("F18").Value = 704
$ws4.Range("F25").Value = 1106
$ws4.Range("F26").Value = 2036
$ws4.Range("F29").Value = 667
$ws4.Range("F31").Value = 5197
$ws4.Range("F33").Value = 1232
$ws4.Range("F34").Value = 3648
$ws4.Range("F37").Value = 1155
$ws4.Range("F41").Value = 953
$ws4.Range("F42").Value = 367
$ws4.Range("F46").Value = 879
$ws4.Range("F47").Value = 1039
